# Helper: convert a target EMU value to the "points" value that, once the
# host round-trips it through a 32-bit COM Single (as real PowerPoint does
# for Shape.Left/Top/Width/Height), truncates back to exactly that EMU
# value instead of landing one EMU short.
function ToPt([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 ("OVERVIEW/GOALS/REQUIREMENTS" -> "Accomplishments")
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$s2Title = $s2.Shapes.Item(1).TextFrame.TextRange
$s2Title.Text = "Accomplishments"
$s2Title.ParagraphFormat.Alignment = 2   # ppAlignCenter
$s2Title.LanguageID = "en-US"

$s2Body = $s2.Shapes.Item(2).TextFrame.TextRange
$s2Body.Paragraphs(1).Text = "Resurrected to previous capabilities to establish baseline."
$s2Body.Paragraphs(2).Text = "Researched software improvements "
$s2Body.Paragraphs(3).Text = "Added Keyboard controls to control robot "
$s2Body.Paragraphs(4).Text = "Added to head and arm movements, more gestures, smoother movement – All working – debug needed in concurrency "
$s2Body.Paragraphs(5).Text = "Connected API’s and webhooks for vocal recognition "
$s2Body.Paragraphs(6).Text = "Improved speech synthesis, and increase library of vocabulary"
$s2Body.Paragraphs(7).Text = "Wrote script to perform in a theater play"

# ---------------------------------------------------------------------
# Slide 3 ("PROJECT TASKS")
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body.Paragraphs(1).Text = "The divided the project tasks into 5 categories:"
$s3Body.Paragraphs(2).Text = "Speech Synthesis – In progress"
$s3Body.Paragraphs(3).Text = "Voice Recognition – In progress – Google Play needs wifi "
$s3Body.Paragraphs(4).Text = "Hardware improvements/repairs – accomplished "
$s3Body.Paragraphs(5).Text = "Robot Theater Scripting – accomplished "
$s3Body.Paragraphs(6).Text = "Director – Keyboard and script control "

# ---------------------------------------------------------------------
# Slide 4 ("SPEECH SYNTHESIS")
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4Body.Paragraphs(3).Text = "Speech synthesis is accomplished using a text to speech web application – Google Play"

# ---------------------------------------------------------------------
# Slide 9 ("FUTURE PLANS")
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9BodyShape = $s9.Shapes.Item(2)

# give the body placeholder an explicit position/size (previously inherited
# from the layout with an empty <p:spPr/>)
$s9BodyShape.Left = ToPt 1368938
$s9BodyShape.Top = ToPt 1238937
$s9BodyShape.Width = ToPt 7038900
$s9BodyShape.Height = ToPt 2911200

$s9Body = $s9BodyShape.TextFrame.TextRange
$s9Body.Paragraphs(4).Text = "Einstein waive his arms and preforms based on commands"
